# Append three new boiler component rows (ID_Boiler 3, 4, 5) below the
# existing data, matching columns: ID_Boiler | type | power_max | power_max_unit | carnot_efficiency_factor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: biomass
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "biomass"
$ws.Range("C4").Value = 15000
$ws.Range("D4").Value = "W"
$ws.Range("E4").Value = 0.35

# Row 5: district_heating
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "district_heating"
$ws.Range("C5").Value = 15000
$ws.Range("D5").Value = "W"
$ws.Range("E5").Value = 0.35

# Row 6: heating_oil
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "heating_oil"
$ws.Range("C6").Value = 15000
$ws.Range("D6").Value = "W"
$ws.Range("E6").Value = 0.35
